# Updates the cryptos list (prices / 1h volume %) with freshly scraped
# figures, matching the "Updated cryptos list ... with GitHub Actions"
# commit. Price cells (column D) hold numeric-looking text (e.g.
# "68.521.10", "0.0000260", "388.10") that must stay plain text, so we
# force NumberFormat to Text before assigning the value and then restore
# the cell to the workbook's default "Normal" style (these cells carry no
# custom style in the source file) so no stray formatting is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextCell "D2" '68.521.10'
$ws.Range("E2").Value = '  +0.61%  '
Set-TextCell "D3" '3.755.98'
$ws.Range("E3").Value = '  -0.72%  '
$ws.Range("E4").Value = '  -0.06%  '
Set-TextCell "D5" '594.13'
$ws.Range("E5").Value = '  -0.62%  '
Set-TextCell "D6" '166.63'
$ws.Range("E6").Value = '  -2.20%  '
Set-TextCell "D7" '3.754.82'
$ws.Range("E7").Value = '  -0.70%  '
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("E9").Value = '  -1.32%  '
$ws.Range("E10").Value = '  -3.01%  '
Set-TextCell "D11" '6.45'
$ws.Range("E11").Value = '  -1.53%  '
$ws.Range("E12").Value = '  -1.38%  '
Set-TextCell "D13" '0.0000260'
$ws.Range("E13").Value = '  -7.73%  '
Set-TextCell "D14" '36.09'
$ws.Range("E14").Value = '  -1.58%  '
Set-TextCell "D15" '4.387.67'
$ws.Range("E15").Value = '  -0.68%  '
Set-TextCell "D16" '3.765.54'
$ws.Range("E16").Value = '  -0.63%  '
Set-TextCell "D17" '68.512.72'
$ws.Range("E17").Value = '  +0.59%  '
$ws.Range("E18").Value = '  -4.83%  '
$ws.Range("E19").Value = '  -0.12%  '
Set-TextCell "D20" '6.98'
$ws.Range("E20").Value = '  -3.39%  '
Set-TextCell "D21" '10.77'
$ws.Range("E21").Value = '  +1.30%  '
Set-TextCell "D22" '464.64'
$ws.Range("E22").Value = '  -0.94%  '
$ws.Range("E23").Value = '  -3.42%  '
Set-TextCell "D24" '84.12'
$ws.Range("E24").Value = '  +0.32%  '
Set-TextCell "D25" '0.0000145'
$ws.Range("E25").Value = '  -3.94%  '
$ws.Range("E26").Value = '  -3.28%  '
Set-TextCell "D27" '11.94'
$ws.Range("E27").Value = '  -1.88%  '
Set-TextCell "D28" '10.02'
$ws.Range("E28").Value = '  -4.90%  '
$ws.Range("E29").Value = '  -0.10%  '
Set-TextCell "D30" '3.904.77'
$ws.Range("E30").Value = '  -0.70%  '
$ws.Range("E31").Value = '  -5.46%  '
Set-TextCell "D32" '7.29'
$ws.Range("E32").Value = '  -4.23%  '
Set-TextCell "D33" '29.97'
$ws.Range("E33").Value = '  -1.94%  '
$ws.Range("E34").Value = '  -3.68%  '
Set-TextCell "D35" '9.17'
$ws.Range("E35").Value = '  -0.72%  '
Set-TextCell "D37" '3.708.83'
$ws.Range("E37").Value = '  -0.99%  '
$ws.Range("E38").Value = '  -4.04%  '
Set-TextCell "D39" '3.38'
$ws.Range("E39").Value = '  -10.24%  '
$ws.Range("E40").Value = '  -0.29%  '
$ws.Range("E41").Value = '  -0.61%  '
$ws.Range("E42").Value = '  -1.27%  '
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("E44").Value = '  -0.01%  '
Set-TextCell "D45" '0.302'
$ws.Range("E45").Value = '  -4.12%  '
Set-TextCell "D46" '43.41'
$ws.Range("E46").Value = '  +7.62%  '
$ws.Range("B47").Value = 'Cosmos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell "D47" '8.49'
$ws.Range("E47").Value = '  -2.46%  '
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextCell "D48" '1.91'
$ws.Range("E48").Value = '  -2.43%  '
Set-TextCell "D49" '46.42'
$ws.Range("E49").Value = '  +1.62%  '
Set-TextCell "D50" '146.02'
$ws.Range("E50").Value = '  +3.74%  '
Set-TextCell "D51" '388.10'
$ws.Range("E51").Value = '  -4.02%  '
